$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report date range) ---
$ws.Range("A8").Value = "Volume 31   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/2/2024  Through  9/8/2024"

# --- Crime statistics table updates ---
$ws.Range("C14").NumberFormat = '#,##0'
$ws.Range("C14").Value = 1
$ws.Range("F14").NumberFormat = '#,##0'
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 50
$ws.Range("L14").Value = -25
$ws.Range("M14").Value = -70
$ws.Range("N14").Value = -78.571428571428
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 14
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 133.333333333333
$ws.Range("N15").Value = -34.375
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 8.333333333333
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 115
$ws.Range("K16").Value = -16.521739130434
$ws.Range("L16").Value = -26.717557251908
$ws.Range("M16").Value = -56.950672645739
$ws.Range("N16").Value = -86.344238975817
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -11.111111111111
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 35
$ws.Range("H17").Value = -31.428571428571
$ws.Range("I17").Value = 298
$ws.Range("J17").Value = 311
$ws.Range("K17").Value = -4.180064308681
$ws.Range("L17").Value = -3.559870550161
$ws.Range("M17").Value = 44.660194174757
$ws.Range("N17").Value = 5.673758865248
$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -85.714285714285
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -59.090909090909
$ws.Range("I18").Value = 122
$ws.Range("J18").Value = 157
$ws.Range("K18").Value = -22.292993630573
$ws.Range("L18").Value = -6.153846153846
$ws.Range("M18").Value = -50.806451612903
$ws.Range("N18").Value = -88.734995383194
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = -50
$ws.Range("I19").Value = 398
$ws.Range("J19").Value = 492
$ws.Range("K19").Value = -19.105691056910
$ws.Range("L19").Value = -10.961968680089
$ws.Range("M19").Value = 31.788079470198
$ws.Range("N19").Value = -3.163017031630
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -44.444444444444
$ws.Range("F20").Value = 43
$ws.Range("G20").Value = 45
$ws.Range("H20").Value = -4.444444444444
$ws.Range("I20").Value = 281
$ws.Range("J20").Value = 238
$ws.Range("K20").Value = 18.067226890756
$ws.Range("L20").Value = 66.272189349112
$ws.Range("M20").Value = 6.439393939393
$ws.Range("N20").Value = -87.893149504523
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -35
$ws.Range("F21").Value = 126
$ws.Range("G21").Value = 184
$ws.Range("H21").Value = -31.521739130434
$ws.Range("I21").Value = 1219
$ws.Range("J21").Value = 1329
$ws.Range("K21").Value = -8.276899924755
$ws.Range("L21").Value = 1.668056713928
$ws.Range("M21").Value = -3.788476716653
$ws.Range("N21").Value = -74.845233182005
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -32.142857142857
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = -40.944881889763
$ws.Range("I24").Value = 766
$ws.Range("J24").Value = 875
$ws.Range("K24").Value = -12.457142857142
$ws.Range("L24").Value = -27.667610953729
$ws.Range("M24").Value = 34.385964912280
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -57.142857142857
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = -42.857142857142
$ws.Range("I25").Value = 181
$ws.Range("J25").Value = 154
$ws.Range("K25").Value = 17.532467532467
$ws.Range("L25").Value = -28.740157480315
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = -31.578947368421
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = 21.568627450980
$ws.Range("I26").Value = 541
$ws.Range("J26").Value = 446
$ws.Range("K26").Value = 21.300448430493
$ws.Range("L26").Value = 19.690265486725
$ws.Range("M26").Value = 22.675736961451
$ws.Range("D14").Copy($ws.Range("C27"))
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("J27").Value = 25
$ws.Range("K27").Value = 36
$ws.Range("L27").Value = 61.904761904761
$ws.Range("D14").Copy($ws.Range("C28"))
$ws.Range("D14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 600
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("C29").Value = 3
$ws.Range("F29").NumberFormat = '#,##0'
$ws.Range("F29").Value = 3
$ws.Range("I29").Value = 14
$ws.Range("K29").Value = 100
$ws.Range("L29").Value = -26.315789473684
$ws.Range("M29").Value = -54.838709677419
$ws.Range("N29").Value = -70.212765957446
$ws.Range("C30").NumberFormat = '#,##0'
$ws.Range("C30").Value = 1
$ws.Range("F30").NumberFormat = '#,##0'
$ws.Range("F30").Value = 1
$ws.Range("I30").Value = 11
$ws.Range("K30").Value = 83.333333333333
$ws.Range("L30").Value = -26.666666666666
$ws.Range("M30").Value = -50
$ws.Range("N30").Value = -74.418604651162
$ws.Range("D14").Copy($ws.Range("C33"))
$ws.Range("D33").NumberFormat = '#,##0'
$ws.Range("D33").Value = 1
$ws.Range("E33").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E33").Value = -100
$ws.Range("J33").Value = 8
$ws.Range("K33").Value = -62.5
